$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update frequency (column C) values for rows 2-9 (category pairs unchanged)
$ws.Range("C2").Value = 3452
$ws.Range("C3").Value = 3092
$ws.Range("C4").Value = 2149
$ws.Range("C5").Value = 1884
$ws.Range("C6").Value = 1342
$ws.Range("C7").Value = 688
$ws.Range("C8").Value = 585
$ws.Range("C9").Value = 497

# Row 10: categories change and frequency changes
$ws.Range("A10").Value = "Textiles & Cozy Items"
$ws.Range("B10").Value = "Textiles & Cozy Items"
$ws.Range("C10").Value = 493

# Row 11: categories change and frequency changes
$ws.Range("A11").Value = "Home Decor"
$ws.Range("B11").Value = "Vintage & Collectibles"
$ws.Range("C11").Value = 447
